# Update "想去人数" (want-to-go count) figures in the 苏州-漫展信息 workbook.
# The same set of events is listed on both the "展览" sheet and the
# consolidated "全部类型" sheet, so the F-column values must be updated in
# both places, even though the row numbers differ between the two sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 2091
$wsExhibit.Range("F5").Value = 367
$wsExhibit.Range("F6").Value = 642
$wsExhibit.Range("F9").Value = 10734
$wsExhibit.Range("F11").Value = 159
$wsExhibit.Range("F15").Value = 7582
$wsExhibit.Range("F17").Value = 725
$wsExhibit.Range("F18").Value = 273
$wsExhibit.Range("F20").Value = 3345

# Sheet "全部类型" (all types) - same events, different row numbers
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 2091
$wsAll.Range("F5").Value = 367
$wsAll.Range("F6").Value = 642
$wsAll.Range("F12").Value = 10734
$wsAll.Range("F14").Value = 159
$wsAll.Range("F18").Value = 7582
$wsAll.Range("F20").Value = 725
$wsAll.Range("F21").Value = 273
$wsAll.Range("F23").Value = 3345
